# Updated list of components with PCB references
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities (column D) for the affected component rows.
# Each of these rows has a formula in column G (Quantity * Unit Cost)
# that will recalculate automatically, which in turn updates the
# grand total formula in B52 (SUM(G2:G64)).
$ws.Range("D21").Value = 5
$ws.Range("D22").Value = 50
$ws.Range("D23").Value = 50
$ws.Range("D24").Value = 50
$ws.Range("D25").Value = 50
$ws.Range("D27").Value = 50
$ws.Range("D28").Value = 50

# Move the active selection to D29 (also clears the previous
# topLeftCell="C1" scroll-freeze state from the saved view).
$ws.Range("D29").Select()
